$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 237, pushing existing rows 237-287 down to 238-288
$ws.Rows.Item(237).Insert()

# Populate the newly inserted row 237 with the new record's data
$ws.Cells.Item(237, 1).Value = 7
$ws.Cells.Item(237, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(237, 3).Value = "Ñuble"
$ws.Cells.Item(237, 4).Value = (Get-Date -Year 2023 -Month 3 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(237, 5).Value = 16
$ws.Cells.Item(237, 6).Value = "Fruta"
$ws.Cells.Item(237, 7).Value = 100108
$ws.Cells.Item(237, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(237, 9).Value = 100108005
$ws.Cells.Item(237, 10).Value = "Piña"
$ws.Cells.Item(237, 11).Value = "Caramelo"
$ws.Cells.Item(237, 12).Value = "Segunda"
$ws.Cells.Item(237, 13).Value = 50
$ws.Cells.Item(237, 14).Value = 23000
$ws.Cells.Item(237, 15).Value = 23000
$ws.Cells.Item(237, 16).Value = 23000
$ws.Cells.Item(237, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(237, 18).Value = "Ecuador"
$ws.Cells.Item(237, 19).Value = 1643
$ws.Cells.Item(237, 20).Value = 14
